# Pedidos.xlsx update — insert/refresh rows 212-216 on the active sheet.
#
# Net effect (per the source diff): a new shipment row (80267042 /
# 10253-ARI-I / 1) is inserted ahead of the existing 84004838 row, and the
# final previously-blank row of the block is filled in with a new
# 84004841 / 10255-ARI-I / 1 row. Concretely this rewrites the A:C values
# for rows 212-216 in place (rows 215-216 were blank placeholder rows that
# already carried the correct cell styles, so nothing below row 216 moves).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a text value into a cell without Excel's COM layer silently
# re-typing a numeric-looking string as a number: assign it as a
# string-literal formula, then freeze the formula down to a static value.
# (Using a "@" text NumberFormat instead would leave a stray unused style
# behind in styles.xml, which we want to avoid.)
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163) # xlPasteValues
}

$rows = @(
    @{ Row = 212; A = "80267042"; B = "10253-ARI-I"; C = 1 },
    @{ Row = 213; A = "84004838"; B = "20953-CTY-I"; C = 3 },
    @{ Row = 214; A = "84004839"; B = "10257-ARI-I"; C = 1 },
    @{ Row = 215; A = "84004840"; B = "10355-ARI-I"; C = 1 },
    @{ Row = 216; A = "84004841"; B = "10255-ARI-I"; C = 1 }
)

foreach ($r in $rows) {
    Set-TextValue $ws.Range("A$($r.Row)") $r.A
    Set-TextValue $ws.Range("B$($r.Row)") $r.B
    $ws.Range("C$($r.Row)").Value = $r.C
}

$excel.CutCopyMode = $false

# Match the saved view state: selection A2:C216 (active cell A2), no
# frozen/scrolled topLeftCell override.
$ws.Activate() | Out-Null
$ws.Range("A2:C216").Select() | Out-Null
